$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the risk-point keyword for the welding row: add "용접기" to the list
$ws.Range("C3").Value = "용접, 슬래그, 용접기"

# Move the active selection to C10, matching the saved view state
$ws.Range("C10").Select()
